$wb = $excel.ActiveWorkbook

# --- Sheet1: add new column R "Ex. Price" = ROUND(F,2) for each trade row ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# Header cell reuses the existing "Ex. Price" shared string (same text as F1)
$ws1.Range("R1").Value = "Ex. Price"

# R2 is a standalone (non-shared) formula
$ws1.Range("R2").Formula = "=ROUND(F2,2)"

# R3:R48 share one formula definition
$ws1.Range("R3:R48").Formula = "=ROUND(F3,2)"

# --- Sheet2 (Tests): a few data corrections ---
$ws2 = $wb.Worksheets.Item("Tests")

# Row 26: exit price rounded down to 2 decimals
$ws2.Range("F26").Value = 56.58

# Row 48: trade type corrected from "Open Short" to "Short" ...
$ws2.Range("B48").Value = "Short"
# ... and the exit date replaced with the text "Open" (still open position)
$ws2.Range("E48").Value = "Open"

# --- View state updates ---
$ws1.Activate()
$ws1.Range("R1").Select()

$ws2.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 26
$win.ScrollColumn = 1
$ws2.Range("E49").Select()

# Sheet1 ends up the active sheet/tab
$ws1.Activate()

Write-Host "done"
